# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# ---------------------------------------------------------------------
# 1) Update the "last updated" timestamp string in A1
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 27 de Junio de 2020 a las 04:41"

# ---------------------------------------------------------------------
# 2) Swap country names that moved position in the ranking
#    (pairs of adjacent rows whose countries were swapped)
# ---------------------------------------------------------------------
$ws.Range("A58").Value  = "Guatemala"
$ws.Range("A59").Value  = "Moldavia"

$ws.Range("A201").Value = "Santa Lucia"
$ws.Range("A202").Value = "Laos"

$ws.Range("A203").Value = "Fiyi"
$ws.Range("A204").Value = "Dominica"

$ws.Range("A208").Value = "Groenlandia"
$ws.Range("A209").Value = "Islas Malvinas"

$ws.Range("A212").Value = "Seychelles"
$ws.Range("A213").Value = "Montserrat"

# ---------------------------------------------------------------------
# 3) Update the statistic figures (Casos totales, Nuevos casos,
#    Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
# ---------------------------------------------------------------------

# Row 47 - Bolivia
$ws.Range("B47").Value = 29423
$ws.Range("C47").Value = 920
$ws.Range("D47").Value = 7736
$ws.Range("E47").Value = 20753
$ws.Range("G47").Value = 21
$ws.Range("H47").Value = 934

# Row 58 - now Guatemala
$ws.Range("B58").Value = 15828
$ws.Range("C58").Value = 209
$ws.Range("D58").Value = 3028
$ws.Range("E58").Value = 12128
$ws.Range("G58").Value = 49
$ws.Range("H58").Value = 672

# Row 59 - now Moldavia
$ws.Range("B59").Value = 15776
$ws.Range("D59").Value = 8765
$ws.Range("E59").Value = 6496
$ws.Range("H59").Value = 515

# Row 65 - Corea del Sur
$ws.Range("B65").Value = 12653
$ws.Range("C65").Value = 51
$ws.Range("D65").Value = 11317
$ws.Range("E65").Value = 1054

# Row 74 - Australia
$ws.Range("B74").Value = 7601
$ws.Range("C74").Value = 6
$ws.Range("D74").Value = 6960
$ws.Range("E74").Value = 537

# Row 81 - Haiti
$ws.Range("D81").Value = 574
$ws.Range("E81").Value = 4873

# Row 212 - now Seychelles
$ws.Range("D212").Value = 11
$ws.Range("H212").Value = 0

# Row 213 - now Montserrat
$ws.Range("D213").Value = 10
$ws.Range("H213").Value = 1
